$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (F1:H1) - Outlier/MAD columns for each algorithm
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Reuse the same header formatting as the existing header row (bold, border, centered)
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)

# New boolean "outlier" flag columns for rows 2-4
$ws.Range("F2").Value = $true
$ws.Range("G2").Value = $false
$ws.Range("H2").Value = $false

$ws.Range("F3").Value = $false
$ws.Range("G3").Value = $true
$ws.Range("H3").Value = $false

$ws.Range("F4").Value = $false
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = $false
